$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$headers = @("资源编号", "资源名称", "资源类型", "资源等级", "资源地图", "资源品质", "专属职业")

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

$ws.Range("I3").Select() | Out-Null
